$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24 and 25 hold the two tumor samples ("7316-2146" and "7316-2751")
# belonging to the same patient (C271092). The row order for this pair is
# being rearranged ("move arrange to end") so the sample that used to sit
# in row 25 now sits in row 24, and vice versa. Columns A, C, D, E, F, G are
# identical between the two rows, so only Tumor ID (B), T/N TelHunt ratio
# (H) and Cohort (I) actually need to swap.

$cols = @("B","H","I")

$row24 = @{}
$row25 = @{}

foreach ($col in $cols) {
    $row24[$col] = $ws.Range("$col`24").Value()
    $row25[$col] = $ws.Range("$col`25").Value()
}

foreach ($col in $cols) {
    $ws.Range("$col`24").Value = $row25[$col]
    $ws.Range("$col`25").Value = $row24[$col]
}
